# Add data for 2022-10-29: the "through" date in the report label/title
# advances from October 20 to October 21, and the month-to-date counts for
# every October on record (2022, 2021, 2020, 2019, 2018, 2017, 2016) are
# refreshed to reflect the new as-of date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the worksheet (tab) name and the report title cell.
$ws.Name = "Through 2022-10-21"
$ws.Range("B1").Value = "October 2022 (through October 21)"

# Garfield Park (row 2)
$ws.Range("V2").Value = 15
$ws.Range("AZ2").Value = 2

# Humboldt Park (row 3)
$ws.Range("V3").Value = 1
$ws.Range("AF3").Value = 3

# North Lawndale (row 5)
$ws.Range("AP5").Value = 2
$ws.Range("AZ5").Value = 4

# Little Village (row 14)
$ws.Range("V14").Value = 3

# Grand Crossing (row 18)
$ws.Range("B18").Value = 2
$ws.Range("V18").Value = 2

# Lower West Side (row 20)
$ws.Range("L20").Value = 3

# Grand Boulevard (row 23)
$ws.Range("L23").Value = 4

# Auburn Gresham (row 24)
$ws.Range("B24").Value = 8

# Lake View (row 26)
$ws.Range("AZ26").Value = 1

# West Town (row 30)
$ws.Range("AZ30").Value = 2
$ws.Range("BJ30").Value = 1

# Portage Park (row 31)
$ws.Range("AF31").Value = 2

# Washington Park (row 36)
$ws.Range("AF36").Value = 1

# Calumet Heights (row 46)
$ws.Range("V46").Value = 2

# Edgewater (row 50)
$ws.Range("AP50").Value = 2

# Mckinley Park (row 53)
$ws.Range("V53").Value = 1

# Chicago Lawn (row 65)
$ws.Range("AZ65").Value = 1

# South Deering (row 94)
$ws.Range("V94").Value = 2
